# Applies the "cleaned and added some inflation rates" edit to Cataloged_Indicators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Remarks for the top Stock Market rows (Stock Indices / Currency Rates) ---
$ws.Range("G2").Value = "Combined Yahoo data with csv from Investing"
$ws.Range("G3").Value = "Combined Yahoo data with csv from Investing"

# --- New GDP rows (16-17): GDP Annual Growth Rate / GDP Growth Rate sourced from
#     OECD, Trading Economics, FX Empire ---
$ws.Range("B16").Value = "GDP Annual Growth Rate"
$ws.Range("C16").Value = "OECD, Trading Economics, FX Empire"
$ws.Range("D16").Value = "Q"
$ws.Range("E16").Value = "Growth YoY"
$ws.Range("F16").Value = 98.4
$ws.Range("G16").Value = "Key indicator, manually filled"

$ws.Range("B17").Value = "GDP Growth Rate"
$ws.Range("C17").Value = "OECD, Trading Economics, FX Empire"
$ws.Range("D17").Value = "Q"
$ws.Range("E17").Value = "Growth QoQ"
$ws.Range("F17").Value = 96.7
$ws.Range("G17").Value = "Key indicator, manually filled"

# --- New Labour rows (23-25): Unemployment growth (IMF) + extra Unemployment Rate ---
$ws.Range("B23").Value = "Unemployment Growth QoQ"
$ws.Range("C23").Value = "IMF"
$ws.Range("D23").Value = "Q"
$ws.Range("E23").Value = "Growth QoQ"
$ws.Range("F23").Value = 84.3

$ws.Range("B24").Value = "Unemployment Growth YoY"
$ws.Range("C24").Value = "IMF"
$ws.Range("D24").Value = "Q"
$ws.Range("E24").Value = "Growth YoY"
$ws.Range("F24").Value = 85.6

$ws.Range("B25").Value = "Unemployment Rate"
$ws.Range("C25").Value = "OECD, Trading Economics, FX Empire"
$ws.Range("D25").Value = "M"
$ws.Range("E25").Value = "Level"
$ws.Range("F25").Value = 94.3
$ws.Range("G25").Value = "Key indicator, manually filled"

# --- New Prices row (37): Inflation Rate sourced from OECD, Trading Economics, FX Empire ---
$ws.Range("B37").Value = "Inflation Rate"
$ws.Range("C37").Value = "OECD, Trading Economics, FX Empire"
$ws.Range("D37").Value = "M"
$ws.Range("E37").Value = "Growth YoY"
$ws.Range("G37").Value = "Key indicator, manually filled"

# --- Clarify the manually-filled PMI country coverage remarks (rows 71-72) ---
$ws.Range("G71").Value = "32 countries from 2012 + 3 from 2013, Manually filled"
$ws.Range("G72").Value = "12 countries from 2012 + 3 from 2013, Manually filled"

# --- Move the active selection to where the author last clicked ---
$ws.Range("F9").Select()
